# Apply hybrid bold + color highlighting to quantitative impact metrics
# (percentages, dollar amounts, large numbers) across achievement and
# work-experience bullet paragraphs.

function Get-WordColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$HighlightColor = Get-WordColor "2C3E50"

function Set-MetricHighlight($Paragraph, [string]$Needle) {
    $rng = $Paragraph.Range
    $found = $rng.Find.Execute($Needle, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $rng.Bold = 1
        $rng.Font.Color = $HighlightColor
    }
    return $found
}

$d = $word.ActiveDocument

# Locate the target bullet paragraphs by matching a distinctive substring,
# so the script is resilient to any paragraph re-numbering.
function Find-ParagraphByText([string]$Snippet) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text.Contains($Snippet)) {
            return $para
        }
    }
    return $null
}

# --- Professional Experience bullets -------------------------------------

$p1 = Find-ParagraphByText "Discovered systematic race coding errors"
Set-MetricHighlight $p1 "23%"  | Out-Null
Set-MetricHighlight $p1 "64%"  | Out-Null

$p2 = Find-ParagraphByText "Utilized advanced sampling methods to decrease survey margin of error from"
Set-MetricHighlight $p2 "±4.2%" | Out-Null
Set-MetricHighlight $p2 "±2.1%" | Out-Null
Set-MetricHighlight $p2 "71%"   | Out-Null
Set-MetricHighlight $p2 "87%"   | Out-Null

$p3 = Find-ParagraphByText "Trigonometric algorithm for boundary estimation"
Set-MetricHighlight $p3 "73.5%" | Out-Null
Set-MetricHighlight $p3 "$4.7M" | Out-Null

$p4 = Find-ParagraphByText "Built real-time FEC analysis systems"
Set-MetricHighlight $p4 "$2"    | Out-Null

# --- Key Achievements and Impact bullets ----------------------------------

$p5 = Find-ParagraphByText "Predictive excellence: Utilized advanced sampling methods"
Set-MetricHighlight $p5 "±4.2%" | Out-Null
Set-MetricHighlight $p5 "±2.1%" | Out-Null

$p6 = Find-ParagraphByText "Increased voter turnout prediction accuracy from"
Set-MetricHighlight $p6 "71%"   | Out-Null
Set-MetricHighlight $p6 "87%"   | Out-Null

$p7 = Find-ParagraphByText "Methodological advancement: Improved segmentation accuracy"
Set-MetricHighlight $p7 "34%"   | Out-Null
Set-MetricHighlight $p7 "28%"   | Out-Null

Write-Host "Metric highlighting applied."
